# Added updated requisites files from Spring 2026 (1261) term
# Sheet1 rows are re-sorted alphabetically by Course_Code and several
# courses get refreshed Course_Name / Prerequisites / CoRequisites / Description values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Course_Code
$ws.Cells.Item(1,1).Value = 'Course_Code'
$ws.Cells.Item(1,2).Value = 'Type'
$ws.Cells.Item(1,3).Value = 'Dept'
$ws.Cells.Item(1,4).Value = 'Course_Name'
$ws.Cells.Item(1,5).Value = 'Prerequisites'
$ws.Cells.Item(1,6).Value = 'CoRequisites'
$ws.Cells.Item(1,7).Value = 'Acad_Group'
$ws.Cells.Item(1,8).Value = 'Acad_Org'
$ws.Cells.Item(1,9).Value = 'Description'

# Row 2: CRIM101
$ws.Cells.Item(2,1).Value = 'CRIM101'
$ws.Cells.Item(2,2).Value = 'Course'
$ws.Cells.Item(2,3).Value = 'CRIM'
$ws.Cells.Item(2,4).Value = 'CRIM 101 - Introduction to Criminology'
$ws.Cells.Item(2,5).Value = 'None'
$ws.Cells.Item(2,6).Value = 'None'
$ws.Cells.Item(2,7).Value = 'ARTS'
$ws.Cells.Item(2,8).Value = 'CRIMINOLGY'
$ws.Cells.Item(2,9).ClearContents()

# Row 3: FALX99
$ws.Cells.Item(3,1).Value = 'FALX99'
$ws.Cells.Item(3,2).Value = 'Course'
$ws.Cells.Item(3,3).Value = 'EDUC'
$ws.Cells.Item(3,4).Value = 'FAL X99 - Foundations of Academic Literacy'
$ws.Cells.Item(3,5).Value = 'None'
$ws.Cells.Item(3,6).Value = 'None'
$ws.Cells.Item(3,7).Value = 'EDUC'
$ws.Cells.Item(3,8).Value = 'EDUCATION'
$ws.Cells.Item(3,9).ClearContents()

# Row 4: INDG101
$ws.Cells.Item(4,1).Value = 'INDG101'
$ws.Cells.Item(4,2).Value = 'Course'
$ws.Cells.Item(4,3).Value = 'INDG'
$ws.Cells.Item(4,4).Value = 'INDG 101 - Introduction to Indigenous Studies'
$ws.Cells.Item(4,5).Value = 'None'
$ws.Cells.Item(4,6).Value = 'None'
$ws.Cells.Item(4,7).Value = 'ARTS'
$ws.Cells.Item(4,8).Value = 'INDIGENOUS'
$ws.Cells.Item(4,9).ClearContents()

# Row 5: INDG102
$ws.Cells.Item(5,1).Value = 'INDG102'
$ws.Cells.Item(5,2).Value = 'Course'
$ws.Cells.Item(5,3).Value = 'INDG'
$ws.Cells.Item(5,4).Value = 'INDG 102 - Indigenous Academic Research'
$ws.Cells.Item(5,5).Value = 'None'
$ws.Cells.Item(5,6).Value = 'None'
$ws.Cells.Item(5,7).Value = 'ARTS'
$ws.Cells.Item(5,8).Value = 'INDIGENOUS'
$ws.Cells.Item(5,9).ClearContents()

# Row 6: INDG110W
$ws.Cells.Item(6,1).Value = 'INDG110W'
$ws.Cells.Item(6,2).Value = 'Course'
$ws.Cells.Item(6,3).Value = 'INDG'
$ws.Cells.Item(6,4).Value = 'INDG 110W - International Indigenous Lifewriting'
$ws.Cells.Item(6,5).Value = 'FALX99'
$ws.Cells.Item(6,6).Value = 'None'
$ws.Cells.Item(6,7).Value = 'ARTS'
$ws.Cells.Item(6,8).Value = 'INDIGENOUS'
$ws.Cells.Item(6,9).ClearContents()

# Row 7: INDG111
$ws.Cells.Item(7,1).Value = 'INDG111'
$ws.Cells.Item(7,2).Value = 'Course'
$ws.Cells.Item(7,3).Value = 'INDG'
$ws.Cells.Item(7,4).Value = 'INDG 111 - Introduction to Indigenous Research Methods'
$ws.Cells.Item(7,5).Value = 'None'
$ws.Cells.Item(7,6).Value = 'None'
$ws.Cells.Item(7,7).Value = 'ARTS'
$ws.Cells.Item(7,8).Value = 'INDIGENOUS'
$ws.Cells.Item(7,9).ClearContents()

# Row 8: INDG201
$ws.Cells.Item(8,1).Value = 'INDG201'
$ws.Cells.Item(8,2).Value = 'Course'
$ws.Cells.Item(8,3).Value = 'INDG'
$ws.Cells.Item(8,4).Value = 'INDG 201 - Canadian Aboriginal Peoples'' Perspectives on History'
$ws.Cells.Item(8,5).Value = 'None'
$ws.Cells.Item(8,6).Value = 'None'
$ws.Cells.Item(8,7).Value = 'ARTS'
$ws.Cells.Item(8,8).Value = 'INDIGENOUS'
$ws.Cells.Item(8,9).ClearContents()

# Row 9: INDG201W
$ws.Cells.Item(9,1).Value = 'INDG201W'
$ws.Cells.Item(9,2).Value = 'Course'
$ws.Cells.Item(9,3).Value = 'INDG'
$ws.Cells.Item(9,4).Value = 'INDG 201W - Indigenous Peoples'' Perspectives on History'
$ws.Cells.Item(9,5).Value = 'FALX99'
$ws.Cells.Item(9,6).Value = 'None'
$ws.Cells.Item(9,7).Value = 'ARTS'
$ws.Cells.Item(9,8).Value = 'INDIGENOUS'
$ws.Cells.Item(9,9).ClearContents()

# Row 10: INDG210
$ws.Cells.Item(10,1).Value = 'INDG210'
$ws.Cells.Item(10,2).Value = 'Course'
$ws.Cells.Item(10,3).Value = 'INDG'
$ws.Cells.Item(10,4).Value = 'INDG 210 - Introduction to Indigenous Digital Media'
$ws.Cells.Item(10,5).Value = 'FALX99,INDG101,INDG201W'
$ws.Cells.Item(10,6).Value = 'INDG101,INDG201W'
$ws.Cells.Item(10,7).Value = 'ARTS'
$ws.Cells.Item(10,8).Value = 'INDIGENOUS'
$ws.Cells.Item(10,9).Value = 'Prerequisite: or Corequisite: INDG101 or INDG201W. Students who took FNST222 ST in Spring 2019 with Dr. Knickerbocker or INDG222 ST in Fall 2020 with Dr. Shield(both titled Introduction to Indigenous Digital Media) cannot take INDG210 for further credit.'

# Row 11: INDG211
$ws.Cells.Item(11,1).Value = 'INDG211'
$ws.Cells.Item(11,2).Value = 'Course'
$ws.Cells.Item(11,3).Value = 'INDG'
$ws.Cells.Item(11,4).Value = 'INDG 211 - Researching Residential Schools: An Analysis of RS in North America'
$ws.Cells.Item(11,5).Value = 'None'
$ws.Cells.Item(11,6).Value = 'None'
$ws.Cells.Item(11,7).Value = 'ARTS'
$ws.Cells.Item(11,8).Value = 'INDIGENOUS'
$ws.Cells.Item(11,9).ClearContents()

# Row 12: INDG212
$ws.Cells.Item(12,1).Value = 'INDG212'
$ws.Cells.Item(12,2).Value = 'Course'
$ws.Cells.Item(12,3).Value = 'INDG'
$ws.Cells.Item(12,4).Value = 'INDG 212 - Indigenous Perceptions of Landscape'
$ws.Cells.Item(12,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(12,6).Value = 'None'
$ws.Cells.Item(12,7).Value = 'ARTS'
$ws.Cells.Item(12,8).Value = 'INDIGENOUS'
$ws.Cells.Item(12,9).ClearContents()

# Row 13: INDG222
$ws.Cells.Item(13,1).Value = 'INDG222'
$ws.Cells.Item(13,2).Value = 'Course'
$ws.Cells.Item(13,3).Value = 'INDG'
$ws.Cells.Item(13,4).Value = 'INDG 222 - Selected Topics in Indigenous Studies'
$ws.Cells.Item(13,5).Value = 'None'
$ws.Cells.Item(13,6).Value = 'None'
$ws.Cells.Item(13,7).Value = 'ARTS'
$ws.Cells.Item(13,8).Value = 'INDIGENOUS'
$ws.Cells.Item(13,9).ClearContents()

# Row 14: INDG232
$ws.Cells.Item(14,1).Value = 'INDG232'
$ws.Cells.Item(14,2).Value = 'Course'
$ws.Cells.Item(14,3).Value = 'INDG'
$ws.Cells.Item(14,4).Value = 'INDG 232 - Indigenous Science'
$ws.Cells.Item(14,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(14,6).Value = 'None'
$ws.Cells.Item(14,7).Value = 'ARTS'
$ws.Cells.Item(14,8).Value = 'INDIGENOUS'
$ws.Cells.Item(14,9).ClearContents()

# Row 15: INDG250
$ws.Cells.Item(15,1).Value = 'INDG250'
$ws.Cells.Item(15,2).Value = 'Course'
$ws.Cells.Item(15,3).Value = 'INDG'
$ws.Cells.Item(15,4).Value = 'INDG 250 - Introduction to Indigenous Policy'
$ws.Cells.Item(15,5).Value = 'INDG101'
$ws.Cells.Item(15,6).Value = 'None'
$ws.Cells.Item(15,7).Value = 'ARTS'
$ws.Cells.Item(15,8).Value = 'INDIGENOUS'
$ws.Cells.Item(15,9).Value = 'REQ: INDG 101.  Students with credit for INDG (or FNST) 222 under the title "Introduction to Public Policy" may not take this course for further credit.'

# Row 16: INDG286
$ws.Cells.Item(16,1).Value = 'INDG286'
$ws.Cells.Item(16,2).Value = 'Course'
$ws.Cells.Item(16,3).Value = 'INDG'
$ws.Cells.Item(16,4).Value = 'INDG 286 - Indigenous Peoples and British Columbia: An Introduction'
$ws.Cells.Item(16,5).Value = 'None'
$ws.Cells.Item(16,6).Value = 'None'
$ws.Cells.Item(16,7).Value = 'ARTS'
$ws.Cells.Item(16,8).Value = 'INDIGENOUS'
$ws.Cells.Item(16,9).ClearContents()

# Row 17: INDG301
$ws.Cells.Item(17,1).Value = 'INDG301'
$ws.Cells.Item(17,2).Value = 'Course'
$ws.Cells.Item(17,3).Value = 'INDG'
$ws.Cells.Item(17,4).Value = 'INDG 301 - Indigenous/Indigenist Research Methods'
$ws.Cells.Item(17,5).Value = 'INDG101,INDG111,INDG201W'
$ws.Cells.Item(17,6).Value = 'None'
$ws.Cells.Item(17,7).Value = 'ARTS'
$ws.Cells.Item(17,8).Value = 'INDIGENOUS'
$ws.Cells.Item(17,9).ClearContents()

# Row 18: INDG305
$ws.Cells.Item(18,1).Value = 'INDG305'
$ws.Cells.Item(18,2).Value = 'Course'
$ws.Cells.Item(18,3).Value = 'INDG'
$ws.Cells.Item(18,4).Value = 'INDG 305 - Treaties in Canada'
$ws.Cells.Item(18,5).Value = 'None'
$ws.Cells.Item(18,6).Value = 'None'
$ws.Cells.Item(18,7).Value = 'ARTS'
$ws.Cells.Item(18,8).Value = 'INDIGENOUS'
$ws.Cells.Item(18,9).ClearContents()

# Row 19: INDG310
$ws.Cells.Item(19,1).Value = 'INDG310'
$ws.Cells.Item(19,2).Value = 'Course'
$ws.Cells.Item(19,3).Value = 'INDG'
$ws.Cells.Item(19,4).Value = 'INDG 310 - Indigenous Film'
$ws.Cells.Item(19,5).Value = 'INDG101'
$ws.Cells.Item(19,6).Value = 'None'
$ws.Cells.Item(19,7).Value = 'ARTS'
$ws.Cells.Item(19,8).Value = 'INDIGENOUS'
$ws.Cells.Item(19,9).Value = 'REQ-Prerequisite: INDG 101.  Students with credit for INDG (or FNST) 322 under the title "Indigenous Film" or "Indigenous Film and Resurgence" may not take this course for further credit.'

# Row 20: INDG322
$ws.Cells.Item(20,1).Value = 'INDG322'
$ws.Cells.Item(20,2).Value = 'Course'
$ws.Cells.Item(20,3).Value = 'INDG'
$ws.Cells.Item(20,4).Value = 'INDG 322 - Special Topics in Indigenous Studies'
$ws.Cells.Item(20,5).Value = 'None'
$ws.Cells.Item(20,6).Value = 'None'
$ws.Cells.Item(20,7).Value = 'ARTS'
$ws.Cells.Item(20,8).Value = 'INDIGENOUS'
$ws.Cells.Item(20,9).ClearContents()

# Row 21: INDG325
$ws.Cells.Item(21,1).Value = 'INDG325'
$ws.Cells.Item(21,2).Value = 'Course'
$ws.Cells.Item(21,3).Value = 'INDG'
$ws.Cells.Item(21,4).Value = 'INDG 325 - History of Indigenous Peoples of North America to 1850'
$ws.Cells.Item(21,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(21,6).Value = 'None'
$ws.Cells.Item(21,7).Value = 'ARTS'
$ws.Cells.Item(21,8).Value = 'INDIGENOUS'
$ws.Cells.Item(21,9).ClearContents()

# Row 22: INDG326
$ws.Cells.Item(22,1).Value = 'INDG326'
$ws.Cells.Item(22,2).Value = 'Course'
$ws.Cells.Item(22,3).Value = 'INDG'
$ws.Cells.Item(22,4).Value = 'INDG 326 - History of Indigenous Peoples of North America Since 1850'
$ws.Cells.Item(22,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(22,6).Value = 'None'
$ws.Cells.Item(22,7).Value = 'ARTS'
$ws.Cells.Item(22,8).Value = 'INDIGENOUS'
$ws.Cells.Item(22,9).ClearContents()

# Row 23: INDG327
$ws.Cells.Item(23,1).Value = 'INDG327'
$ws.Cells.Item(23,2).Value = 'Course'
$ws.Cells.Item(23,3).Value = 'INDG'
$ws.Cells.Item(23,4).Value = 'INDG 327 - Indigenous Women in Canada'
$ws.Cells.Item(23,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(23,6).Value = 'None'
$ws.Cells.Item(23,7).Value = 'ARTS'
$ws.Cells.Item(23,8).Value = 'INDIGENOUS'
$ws.Cells.Item(23,9).Value = 'REQ-45 units and one of INDG 101 or 201W, or permission of instructor. Students with INDG 322 under this topic may not take this course for further credit. INDG 327 and GSWS 327 are identical and students may not take both courses for credit.'

# Row 24: INDG329
$ws.Cells.Item(24,1).Value = 'INDG329'
$ws.Cells.Item(24,2).Value = 'Course'
$ws.Cells.Item(24,3).Value = 'INDG'
$ws.Cells.Item(24,4).Value = 'INDG 329 - Sexuality and Gender: Indigenous Perspectives'
$ws.Cells.Item(24,5).Value = 'None'
$ws.Cells.Item(24,6).Value = 'None'
$ws.Cells.Item(24,7).Value = 'ARTS'
$ws.Cells.Item(24,8).Value = 'INDIGENOUS'
$ws.Cells.Item(24,9).ClearContents()

# Row 25: INDG332
$ws.Cells.Item(25,1).Value = 'INDG332'
$ws.Cells.Item(25,2).Value = 'Course'
$ws.Cells.Item(25,3).Value = 'INDG'
$ws.Cells.Item(25,4).Value = 'INDG 332 - Indigenous Ethnobotany'
$ws.Cells.Item(25,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(25,6).Value = 'None'
$ws.Cells.Item(25,7).Value = 'ARTS'
$ws.Cells.Item(25,8).Value = 'INDIGENOUS'
$ws.Cells.Item(25,9).ClearContents()

# Row 26: INDG333
$ws.Cells.Item(26,1).Value = 'INDG333'
$ws.Cells.Item(26,2).Value = 'Course'
$ws.Cells.Item(26,3).Value = 'INDG'
$ws.Cells.Item(26,4).Value = 'INDG 333 - Indigenous Ethnozoology'
$ws.Cells.Item(26,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(26,6).Value = 'None'
$ws.Cells.Item(26,7).Value = 'ARTS'
$ws.Cells.Item(26,8).Value = 'INDIGENOUS'
$ws.Cells.Item(26,9).ClearContents()

# Row 27: INDG353W
$ws.Cells.Item(27,1).Value = 'INDG353W'
$ws.Cells.Item(27,2).Value = 'Course'
$ws.Cells.Item(27,3).Value = 'INDG'
$ws.Cells.Item(27,4).Value = 'INDG 353W - Indigenous Heritage Stewardship'
$ws.Cells.Item(27,5).Value = 'FALX99'
$ws.Cells.Item(27,6).Value = 'None'
$ws.Cells.Item(27,7).Value = 'ARTS'
$ws.Cells.Item(27,8).Value = 'INDIGENOUS'
$ws.Cells.Item(27,9).ClearContents()

# Row 28: INDG360
$ws.Cells.Item(28,1).Value = 'INDG360'
$ws.Cells.Item(28,2).Value = 'Course'
$ws.Cells.Item(28,3).Value = 'INDG'
$ws.Cells.Item(28,4).Value = 'INDG 360 - Popular Writing by Indigenous Authors'
$ws.Cells.Item(28,5).Value = 'None'
$ws.Cells.Item(28,6).Value = 'None'
$ws.Cells.Item(28,7).Value = 'ARTS'
$ws.Cells.Item(28,8).Value = 'INDIGENOUS'
$ws.Cells.Item(28,9).ClearContents()

# Row 29: INDG363
$ws.Cells.Item(29,1).Value = 'INDG363'
$ws.Cells.Item(29,2).Value = 'Course'
$ws.Cells.Item(29,3).Value = 'INDG'
$ws.Cells.Item(29,4).Value = 'INDG 363 - Indigenous Poetry, Poetics, Printmaking'
$ws.Cells.Item(29,5).Value = 'None'
$ws.Cells.Item(29,6).Value = 'None'
$ws.Cells.Item(29,7).Value = 'ARTS'
$ws.Cells.Item(29,8).Value = 'INDIGENOUS'
$ws.Cells.Item(29,9).ClearContents()

# Row 30: INDG383
$ws.Cells.Item(30,1).Value = 'INDG383'
$ws.Cells.Item(30,2).Value = 'Course'
$ws.Cells.Item(30,3).Value = 'INDG'
$ws.Cells.Item(30,4).Value = 'INDG 383 - Indigenous Technology: Art and Sustainability'
$ws.Cells.Item(30,5).Value = 'None'
$ws.Cells.Item(30,6).Value = 'None'
$ws.Cells.Item(30,7).Value = 'ARTS'
$ws.Cells.Item(30,8).Value = 'INDIGENOUS'
$ws.Cells.Item(30,9).Value = 'REQ-45 units or permission of instructor. No prior artistic training/experience required. Students with credit for INDG/FNST 322 (topic: Indigenous Expressive Arts - crafts focus) or FNST 383 may not take for further credit.'

# Row 31: INDG401
$ws.Cells.Item(31,1).Value = 'INDG401'
$ws.Cells.Item(31,2).Value = 'Course'
$ws.Cells.Item(31,3).Value = 'INDG'
$ws.Cells.Item(31,4).Value = 'INDG 401 - Indigenous Peoples and Public Policy'
$ws.Cells.Item(31,5).Value = 'INDG101,INDG201W,INDG250'
$ws.Cells.Item(31,6).Value = 'None'
$ws.Cells.Item(31,7).Value = 'ARTS'
$ws.Cells.Item(31,8).Value = 'INDIGENOUS'
$ws.Cells.Item(31,9).Value = 'REQ-45 units, INDG (or FNST) 101 and one of INDG (or FNST) 201W or 250, or permission of the instructor.  Students with credit for FNST 401 may not take this course for further credit.'

# Row 32: INDG402W
$ws.Cells.Item(32,1).Value = 'INDG402W'
$ws.Cells.Item(32,2).Value = 'Course'
$ws.Cells.Item(32,3).Value = 'INDG'
$ws.Cells.Item(32,4).Value = 'INDG 402W - Discourses of Indigenous Peoples'
$ws.Cells.Item(32,5).Value = 'FALX99,INDG101,INDG201W'
$ws.Cells.Item(32,6).Value = 'None'
$ws.Cells.Item(32,7).Value = 'ARTS'
$ws.Cells.Item(32,8).Value = 'INDIGENOUS'
$ws.Cells.Item(32,9).Value = 'REQ-60 units and one of INDG (or FNST) 101 or 201W, or permission of the instructor. Students with credit for INDG (or FNST) 402 or FNST 402W may not take this course for further credit.'

# Row 33: INDG403
$ws.Cells.Item(33,1).Value = 'INDG403'
$ws.Cells.Item(33,2).Value = 'Course'
$ws.Cells.Item(33,3).Value = 'INDG'
$ws.Cells.Item(33,4).Value = 'INDG 403 - Indigenous Knowledges'
$ws.Cells.Item(33,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(33,6).Value = 'None'
$ws.Cells.Item(33,7).Value = 'ARTS'
$ws.Cells.Item(33,8).Value = 'INDIGENOUS'
$ws.Cells.Item(33,9).ClearContents()

# Row 34: INDG410
$ws.Cells.Item(34,1).Value = 'INDG410'
$ws.Cells.Item(34,2).Value = 'Course'
$ws.Cells.Item(34,3).Value = 'INDG'
$ws.Cells.Item(34,4).Value = 'INDG 410 - Elements of Indigenous Style: Indigenous Editing Practices'
$ws.Cells.Item(34,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(34,6).Value = 'INDG101,INDG201W'
$ws.Cells.Item(34,7).Value = 'ARTS'
$ws.Cells.Item(34,8).Value = 'INDIGENOUS'
$ws.Cells.Item(34,9).Value = 'REQ: or Corequisite: INDG 101 or 201W. Students with credit for PUB 410 or PUB 480 under the title "Indigenous Editing" offered in Spring 2022 may not take this course for further credit.'

# Row 35: INDG419
$ws.Cells.Item(35,1).Value = 'INDG419'
$ws.Cells.Item(35,2).Value = 'Course'
$ws.Cells.Item(35,3).Value = 'INDG'
$ws.Cells.Item(35,4).Value = 'INDG 419 - Aboriginal/Indigenous Justice'
$ws.Cells.Item(35,5).Value = 'CRIM101,INDG101,INDG201W'
$ws.Cells.Item(35,6).Value = 'None'
$ws.Cells.Item(35,7).Value = 'ARTS'
$ws.Cells.Item(35,8).Value = 'INDIGENOUS'
$ws.Cells.Item(35,9).ClearContents()

# Row 36: INDG422
$ws.Cells.Item(36,1).Value = 'INDG422'
$ws.Cells.Item(36,2).Value = 'Course'
$ws.Cells.Item(36,3).Value = 'INDG'
$ws.Cells.Item(36,4).Value = 'INDG 422 - Special Topics in Indigenous Studies'
$ws.Cells.Item(36,5).Value = 'None'
$ws.Cells.Item(36,6).Value = 'None'
$ws.Cells.Item(36,7).Value = 'ARTS'
$ws.Cells.Item(36,8).Value = 'INDIGENOUS'
$ws.Cells.Item(36,9).ClearContents()

# Row 37: INDG429
$ws.Cells.Item(37,1).Value = 'INDG429'
$ws.Cells.Item(37,2).Value = 'Course'
$ws.Cells.Item(37,3).Value = 'INDG'
$ws.Cells.Item(37,4).Value = 'INDG 429 - Indigenous Peoples and International Law'
$ws.Cells.Item(37,5).Value = 'CRIM101,INDG101,INDG201W'
$ws.Cells.Item(37,6).Value = 'None'
$ws.Cells.Item(37,7).Value = 'ARTS'
$ws.Cells.Item(37,8).Value = 'INDIGENOUS'
$ws.Cells.Item(37,9).ClearContents()

# Row 38: INDG433
$ws.Cells.Item(38,1).Value = 'INDG433'
$ws.Cells.Item(38,2).Value = 'Course'
$ws.Cells.Item(38,3).Value = 'INDG'
$ws.Cells.Item(38,4).Value = 'INDG 433 - Indigenous Environmental Justice and Activism'
$ws.Cells.Item(38,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(38,6).Value = 'None'
$ws.Cells.Item(38,7).Value = 'ARTS'
$ws.Cells.Item(38,8).Value = 'INDIGENOUS'
$ws.Cells.Item(38,9).Value = 'REQ-45 units and one of INDG (or FNST) 101 or 201W, or permission of the instructor.  Students with credit for FNST 433 may not take this course for further credit.'

# Row 39: INDG435
$ws.Cells.Item(39,1).Value = 'INDG435'
$ws.Cells.Item(39,2).Value = 'Course'
$ws.Cells.Item(39,3).Value = 'INDG'
$ws.Cells.Item(39,4).Value = 'INDG 435 - STT-Land-Based Learning Field Course'
$ws.Cells.Item(39,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(39,6).Value = 'None'
$ws.Cells.Item(39,7).Value = 'ARTS'
$ws.Cells.Item(39,8).Value = 'INDIGENOUS'
$ws.Cells.Item(39,9).Value = 'REQ-INDG 101 or INDG 201W, and permission of the instructor. Students must successfully complete a Criminal Record Check.'

# Row 40: INDG442
$ws.Cells.Item(40,1).Value = 'INDG442'
$ws.Cells.Item(40,2).Value = 'Course'
$ws.Cells.Item(40,3).Value = 'INDG'
$ws.Cells.Item(40,4).Value = 'INDG 442 - Directed Readings in Indigenous Studies'
$ws.Cells.Item(40,5).Value = 'None'
$ws.Cells.Item(40,6).Value = 'None'
$ws.Cells.Item(40,7).Value = 'ARTS'
$ws.Cells.Item(40,8).Value = 'INDIGENOUS'
$ws.Cells.Item(40,9).ClearContents()

# Row 41: INDG443W
$ws.Cells.Item(41,1).Value = 'INDG443W'
$ws.Cells.Item(41,2).Value = 'Course'
$ws.Cells.Item(41,3).Value = 'INDG'
$ws.Cells.Item(41,4).Value = 'INDG 443W - Aboriginal Peoples, History and the Law'
$ws.Cells.Item(41,5).Value = 'FALX99,INDG101,INDG201W'
$ws.Cells.Item(41,6).Value = 'None'
$ws.Cells.Item(41,7).Value = 'ARTS'
$ws.Cells.Item(41,8).Value = 'INDIGENOUS'
$ws.Cells.Item(41,9).ClearContents()

# Row 42: INDG447
$ws.Cells.Item(42,1).Value = 'INDG447'
$ws.Cells.Item(42,2).Value = 'Course'
$ws.Cells.Item(42,3).Value = 'INDG'
$ws.Cells.Item(42,4).Value = 'INDG 447 - Directed Studies in Indigenous Studies'
$ws.Cells.Item(42,5).Value = 'None'
$ws.Cells.Item(42,6).Value = 'None'
$ws.Cells.Item(42,7).Value = 'ARTS'
$ws.Cells.Item(42,8).Value = 'INDIGENOUS'
$ws.Cells.Item(42,9).ClearContents()

# Row 43: INDG462
$ws.Cells.Item(43,1).Value = 'INDG462'
$ws.Cells.Item(43,2).Value = 'Course'
$ws.Cells.Item(43,3).Value = 'INDG'
$ws.Cells.Item(43,4).Value = 'INDG 462 - Indigenous Oral Testimony: Theory, Practice, Purpose, Community'
$ws.Cells.Item(43,5).Value = 'INDG101,INDG201W'
$ws.Cells.Item(43,6).Value = 'None'
$ws.Cells.Item(43,7).Value = 'ARTS'
$ws.Cells.Item(43,8).Value = 'INDIGENOUS'
$ws.Cells.Item(43,9).Value = 'REQ-60 units including INDG (or FNST) 101 or INDG (or FNST) 201W or permission of the instructor.  Students with credit for FNST 462 may not take this course for further credit.'

# Update the sheet selection left by the author after their last edit
$ws.Range("A2:I43").Select()
